# Generate Report for Handback
# Updates the "Latest Handoff Datetime" / "Latest Handback DateTime" timestamps
# recorded for the 66a28a16-a372-45bd-affb-1fc4add6ce51 file across the
# Overview, zh-cn and de-de worksheets to reflect the freshly generated report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-18 18:50:45"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-18 18:50:39"
$zhcn.Range("K2").Value = "2016-08-18 18:50:57"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-18 18:50:45"
$dede.Range("K2").Value = "2016-08-18 18:51:13"
